# Updated symbol list on Tue Feb 14 23:33:23 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'296.25"
$ws.Range('E2').Value = "'1.29%"

$ws.Range('D3').Value = "'42.35"
$ws.Range('E3').Value = "'3.91%"

$ws.Range('D4').Value = "'5.036"
$ws.Range('E4').Value = "'0.02%"

$ws.Range('D5').Value = "'0.07589"
$ws.Range('E5').Value = "'2.63%"

$ws.Range('B6').Value = 'GateToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D6').Value = "'4.393"
$ws.Range('E6').Value = "'2.63%"

$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D7').Value = "'1.612"
$ws.Range('E7').Value = "'3.94%"

$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').Value = "'0.9314"
$ws.Range('E8').Value = "'0.67%"

$ws.Range('B9').Value = 'BTSEToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D9').Value = "'2.405"
$ws.Range('E9').Value = "'3.31%"

$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = "'0.1210"
$ws.Range('E10').Value = "'4.96%"

$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = "'0.1841"
$ws.Range('E11').Value = "'6.61%"

$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = "'0.08979"
$ws.Range('E12').Value = "'3.64%"

$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = "'0.04000"
$ws.Range('E13').Value = "'-4.10%"

$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = "'0.1052"
$ws.Range('E14').Value = "'-0.34%"

$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = "'0.001288"
$ws.Range('E15').Value = "'1.13%"

$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = "'0.005788"
$ws.Range('E16').Value = "'-2.20%"

$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = "'3.366"
$ws.Range('E17').Value = "'-1.41%"

$ws.Range('E18').Value = "'1.11%"

$ws.Range('D19').Value = "'7.870"
$ws.Range('E19').Value = "'2.12%"

$ws.Range('D20').Value = "'0.1420"
$ws.Range('E20').Value = "'3.06%"

$ws.Range('D21').Value = "'0.3001"
$ws.Range('E21').Value = "'4.36%"

$ws.Range('D22').Value = "'0.04054"
$ws.Range('E22').Value = "'5.06%"

$ws.Range('D23').Value = "'0.001267"
$ws.Range('E23').Value = "'0.79%"

$ws.Range('D24').Value = "'0.003981"
$ws.Range('E24').Value = "'3.25%"

$ws.Range('E25').Value = "'-3.70%"

$ws.Range('E26').Value = "'0.26%"

$ws.Range('D38').Value = "'0.02422"
$ws.Range('E38').Value = "'3.29%"

$ws.Range('D39').Value = "'0.05213"
$ws.Range('E39').Value = "'3.93%"

$ws.Range('D40').Value = "'0.006061"
$ws.Range('E40').Value = "'5.07%"

$ws.Range('D41').Value = "'0.007758"
$ws.Range('E41').Value = "'1.15%"

$ws.Range('E42').Value = "'3.71%"

$ws.Range('D43').Value = "'0.007537"
$ws.Range('E43').Value = "'2.87%"

$ws.Range('D44').Value = "'0.007251"
$ws.Range('E44').Value = "'2.15%"

$ws.Range('D45').Value = "'0.2969"
$ws.Range('E45').Value = "'-5.98%"

$ws.Range('D46').Value = "'0.00006788"
$ws.Range('E46').Value = "'5.93%"

$ws.Range('E47').Value = "'0.22%"

$ws.Range('D48').Value = "'0.04611"
$ws.Range('E48').Value = "'175.31%"

$ws.Range('D49').Value = "'0.004202"
$ws.Range('E49').Value = "'0.02%"

$ws.Range('D50').Value = "'0.00002100"
$ws.Range('E50').Value = "'0.22%"

$ws.Range('D51').Value = "'0.0002000"
$ws.Range('E51').Value = "'0.22%"
